# Clean up the "Authors" column (column E) on Sheet1, rows 2-11.
#
# The stored author lists separate entries with a comma followed by a run
# of spaces (e.g. "Name1%..., <spaces>Name2%...."). The fix inserts one
# extra space right after every comma in each of these cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 5)   # column E = Authors
    $old = $cell.Value()
    $new = [System.Text.RegularExpressions.Regex]::Replace($old, ",", ", ")
    $cell.Value = $new
}
